$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.351.58'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.718.26'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4728'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2637'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06210'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.716.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07070'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('E12').Value = '  +3.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5916'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.412'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.347.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006812'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.935.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.550'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.761'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.318'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.405'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '108.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.754'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.003'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07742'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04433'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.615'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9795'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6191'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9338'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '113.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +16.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.415'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.914'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01475'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.327'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.07%  '
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1169'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.290'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05285'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.695'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.217'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('E51').Value = '  +0.80%  '
